$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Technologies")
$ws.Range("M1").Value = "CCS"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").Font.ThemeColor = 1
$ws.Range("M1").Interior.ThemeColor = 9
